# Insert a new "Match ID" column at the front of the passing table.
# (Mirrors: select column A, Insert Shift-Right, fill header + match id
# value of 22 for every data row, bold the new column like the other
# header/id columns, leave the hidden totals row (20) with default style.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$null = $ws.Columns("A").Insert()

# Header cell for the new column.
$ws.Range("A3").Value = "Match ID"

# Bold the header + visible data rows (A3:A19) to match the existing
# "Player ID" style used one column over.
$ws.Range("A3:A19").Font.Bold = $true

# Every visible player row shares the same match id.
$ws.Range("A4:A19").Value = 22

# Hidden totals row keeps the default (unbolded) style.
$ws.Range("A20").Value = 22
$null = $ws.Rows("20").AutoFit()

# Restore the selection to the newly inserted column's data cells.
$null = $ws.Range("A3:A19").Select()
